$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the lab 4 topic text (row 6, column C) to the markdown link version.
$ws.Range("C6").Value = "[Lab 4 Distributions II](https://crumplab.github.io/psyc7709Lab/articles/Lab4_Distributions_II.html)"

# The wrapped link text now needs two lines, so the row grows to match
# (mirrors row 2, which holds a similar markdown-link entry at ht=34).
$ws.Rows.Item(6).RowHeight = 34

# Move the selection to C7, matching the post-edit cursor position.
$ws.Range("C7").Select()
